# Update "想去人数" (interest count) figures on each sheet to the newly
# scraped values (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value  = 14019
$ws.Range("F10").Value = 566
$ws.Range("F12").Value = 13
$ws.Range("F14").Value = 14318
$ws.Range("F15").Value = 385
$ws.Range("F17").Value = 15082
$ws.Range("F19").Value = 8449
$ws.Range("F30").Value = 1056
$ws.Range("F35").Value = 26
$ws.Range("F40").Value = 238
$ws.Range("F41").Value = 400
$ws.Range("F43").Value = 5208

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 54

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value  = 14019
$ws.Range("F10").Value = 566
$ws.Range("F12").Value = 13
$ws.Range("F14").Value = 14318
$ws.Range("F15").Value = 385
$ws.Range("F17").Value = 15082
$ws.Range("F19").Value = 8449
$ws.Range("F31").Value = 1056
$ws.Range("F36").Value = 26
$ws.Range("F37").Value = 54
$ws.Range("F43").Value = 238
$ws.Range("F44").Value = 400
$ws.Range("F46").Value = 5208
